# Sablon.word_ml no longer inserts nested paragraphs.
#
# The template paragraph used to hold a single run whose <w:rPr> carried
# <w:noProof/> while the run body itself smuggled two whole <w:p> elements
# (plus a <w:sectPr>) inside it -- that's how the old Sablon word_ml
# inserted "I *am* / a **Software Developer**". We now flatten that to
# plain sibling paragraphs, and the stray <w:sectPr> that used to live
# inside the nested content collapses away -- its <w:pgNumType w:start="7"/>
# needs to move onto the document's real (last) section properties.

$d = $word.ActiveDocument

# 1. Replace the paragraph that contains the nested-paragraph WordML blob
#    with the flattened WordML: two ordinary sibling paragraphs, no nested
#    <w:p>, no <w:noProof/>, no stray <w:sectPr>.
$p = $d.Paragraphs.Item(5)
$wordMl = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
            '<w:r><w:t xml:space="preserve">I </w:t></w:r>' + `
            '<w:r><w:rPr><w:i/></w:rPr><w:t>am</w:t></w:r>' + `
            '<w:r><w:t xml:space="preserve"> </w:t></w:r>' + `
          '</w:p>' + `
          '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
            '<w:r><w:t xml:space="preserve">a </w:t></w:r>' + `
            '<w:r><w:rPr><w:b/></w:rPr><w:t>Software Developer</w:t></w:r>' + `
          '</w:p>'
$p.Range.InsertXML($wordMl)

# 2. The removed nested <w:sectPr> used to restart page numbering at 7;
#    carry that over to the document's actual (last) section.
$section = $d.Sections.First
$header = $section.Headers.Item(1)
$header.PageNumbers.StartingNumber = 7

# 3. The header's "Generated at" timestamp was refreshed.
$header.Range.Find.Execute("10.04.2015 10:03", $true, $false, $false, $false, `
                            $false, $true, 1, $false, "10.04.2015 11:43", 2)
